$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a single new row at row 21, pushing the "motor" section (old rows 21-25)
# down to rows 22-26, preserving their values/styles automatically.
$ws.Rows("21:21").Insert()

# New "R" section occupies the freshly available rows 18-20 (between the "Neo"
# block and the blank separator row that is now row 21).
$ws.Range("A18").Value2 = "R, aka the rotary sensors"
$ws.Range("A19").Value2 = "R1"
$ws.Range("B19").Value2 = 4
$ws.Range("A20").Value2 = "R2"
$ws.Range("B20").Value2 = 12

# Style the new header cell A18: white (theme lt1) font on blue (theme accent1) fill,
# matching the style used for the other section headers (e.g. "Neo").
$ws.Range("A18").Font.ThemeColor = 2
$ws.Range("A18").Interior.ThemeColor = 5

# Page setup tweak (paper size -> Letter) that accompanied this edit
$ws.PageSetup.PaperSize = 9

$ws.Range("B20").Select()
